# Auto-generated edit script applying the Ifrit_Profits market-data refresh
# across the ALC, ARM, BSM, CUL, LTW and WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4423.077
$ws.Range("H80").Value = 269.83334
$ws.Range("I80").Value = 235.4
$ws.Range("J80").Value = 312.875
$ws.Range("K80").Value = 706.2
$ws.Range("L80").Value = 938.625
$ws.Range("M80").Value = 291.8
$ws.Range("N80").Value = -2934.625
$ws.Range("H83").Value = 269.83334
$ws.Range("I83").Value = 235.4
$ws.Range("J83").Value = 312.875
$ws.Range("K83").Value = 2118.6
$ws.Range("L83").Value = 2815.875
$ws.Range("M83").Value = 2873.4
$ws.Range("N83").Value = -12799.875
$ws.Range("H100").Value = 10205832
$ws.Range("I100").Value = 14287354
$ws.Range("J100").Value = 2025
$ws.Range("K100").Value = 14287354
$ws.Range("L100").Value = 2025
$ws.Range("M100").Value = -14286813
$ws.Range("N100").Value = -3107
$ws.Range("H138").Value = 1494741.8
$ws.Range("J138").Value = 4043.1
$ws.Range("L138").Value = 12129.3
$ws.Range("N138").Value = -22409.3
$ws.Range("H141").Value = 820.94446
$ws.Range("I141").Value = 669.4666999999999
$ws.Range("K141").Value = 2008.4001
$ws.Range("M141").Value = 3171.5999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 9800
$ws.Range("J12").Value = 9800
$ws.Range("L12").Value = 9800
$ws.Range("N12").Value = -10146
$ws.Range("H32").Value = 6771.886
$ws.Range("I32").Value = 6924.075
$ws.Range("J32").Value = 5250
$ws.Range("K32").Value = 6924.075
$ws.Range("L32").Value = 5250
$ws.Range("M32").Value = -6637.075
$ws.Range("N32").Value = -5824
$ws.Range("H74").Value = 6902205
$ws.Range("I74").Value = 13334182
$ws.Range("J74").Value = 10800.929
$ws.Range("K74").Value = 13334182
$ws.Range("L74").Value = 10800.929
$ws.Range("M74").Value = -13333308
$ws.Range("N74").Value = -12548.929
$ws.Range("H77").Value = 6902205
$ws.Range("I77").Value = 13334182
$ws.Range("J77").Value = 10800.929
$ws.Range("K77").Value = 66670910
$ws.Range("L77").Value = 54004.645
$ws.Range("M77").Value = -66666542
$ws.Range("N77").Value = -62740.645
$ws.Range("H122").Value = 1526.5
$ws.Range("I122").Value = 1442.4
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 4327.200000000001
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -1877.200000000001
$ws.Range("N122").Value = -9899.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24748.512
$ws.Range("I134").Value = 32856
$ws.Range("J134").Value = 1163.091
$ws.Range("K134").Value = 98568
$ws.Range("L134").Value = 3489.273
$ws.Range("M134").Value = -96033
$ws.Range("N134").Value = -8559.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 8850.25
$ws.Range("I16").Value = 8850.25
$ws.Range("K16").Value = 26550.75
$ws.Range("M16").Value = -26377.75
$ws.Range("H22").Value = 2416.6667
$ws.Range("I22").Value = 5250
$ws.Range("K22").Value = 15750
$ws.Range("M22").Value = -15581
$ws.Range("H27").Value = 2416.6667
$ws.Range("I27").Value = 5250
$ws.Range("K27").Value = 15750
$ws.Range("M27").Value = -15648
$ws.Range("H31").Value = 1307.1428
$ws.Range("I31").Value = 950
$ws.Range("J31").Value = 2200
$ws.Range("K31").Value = 2850
$ws.Range("L31").Value = 6600
$ws.Range("M31").Value = -2562
$ws.Range("N31").Value = -7176
$ws.Range("H35").Value = 1896.6666
$ws.Range("J35").Value = 1896.6666
$ws.Range("L35").Value = 5689.9998
$ws.Range("N35").Value = -6265.9998
$ws.Range("H74").Value = 2250
$ws.Range("I74").Value = 1500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3439
$ws.Range("H77").Value = 2250
$ws.Range("I77").Value = 1500
$ws.Range("K77").Value = 13500
$ws.Range("M77").Value = -8196
$ws.Range("H93").Value = 41499.5
$ws.Range("J93").Value = 2999
$ws.Range("L93").Value = 8997
$ws.Range("N93").Value = -12741
$ws.Range("H94").Value = 4427
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4427
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 13281
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -14633
$ws.Range("H95").Value = 3300
$ws.Range("J95").Value = 3300
$ws.Range("L95").Value = 9900
$ws.Range("N95").Value = -14018
$ws.Range("H96").Value = 70707200
$ws.Range("J96").Value = 70707200
$ws.Range("L96").Value = 212121600
$ws.Range("N96").Value = -212125718
$ws.Range("H101").Value = 6100.857
$ws.Range("J101").Value = 6913.3335
$ws.Range("L101").Value = 20740.0005
$ws.Range("N101").Value = -25608.0005
$ws.Range("H102").Value = 6333.3335
$ws.Range("J102").Value = 6333.3335
$ws.Range("L102").Value = 19000.0005
$ws.Range("N102").Value = -23868.0005
$ws.Range("H103").Value = 4857763.5
$ws.Range("I103").Value = 5667224
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 17001672
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -17000793
$ws.Range("N103").Value = -4758
$ws.Range("H105").Value = 151335840
$ws.Range("J105").Value = 151335840
$ws.Range("L105").Value = 454007520
$ws.Range("N105").Value = -454012762
$ws.Range("H106").Value = 11115511
$ws.Range("J106").Value = 11115511
$ws.Range("L106").Value = 33346533
$ws.Range("N106").Value = -33348425
$ws.Range("H108").Value = 1249.4
$ws.Range("I108").Value = 1061.75
$ws.Range("J108").Value = 2000
$ws.Range("K108").Value = 3185.25
$ws.Range("L108").Value = 6000
$ws.Range("M108").Value = -305.25
$ws.Range("N108").Value = -11760
$ws.Range("H110").Value = 3660
$ws.Range("J110").Value = 3660
$ws.Range("L110").Value = 10980
$ws.Range("N110").Value = -19160
$ws.Range("H112").Value = 47620284
$ws.Range("I112").Value = 1162.5
$ws.Range("J112").Value = 111112450
$ws.Range("K112").Value = 3487.5
$ws.Range("L112").Value = 333337350
$ws.Range("M112").Value = -2379.5
$ws.Range("N112").Value = -333339566
$ws.Range("H113").Value = 494.28333
$ws.Range("I113").Value = 470.54544
$ws.Range("J113").Value = 523.2963
$ws.Range("K113").Value = 1411.63632
$ws.Range("L113").Value = 1569.8889
$ws.Range("M113").Value = 758.3636799999999
$ws.Range("N113").Value = -5909.8889
$ws.Range("H116").Value = 1875
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 1875
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 5625
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -12509
$ws.Range("H118").Value = 1268
$ws.Range("I118").Value = 741.6
$ws.Range("K118").Value = 2224.8
$ws.Range("M118").Value = -981.8000000000002
$ws.Range("H119").Value = 1489.6
$ws.Range("I119").Value = 1489.6
$ws.Range("K119").Value = 4468.799999999999
$ws.Range("M119").Value = 369.2000000000007
$ws.Range("H120").Value = 10499.111
$ws.Range("I120").Value = 3742.9
$ws.Range("J120").Value = 18944.375
$ws.Range("K120").Value = 11228.7
$ws.Range("L120").Value = 56833.125
$ws.Range("M120").Value = -6390.700000000001
$ws.Range("N120").Value = -66509.125
$ws.Range("H122").Value = 52085316
$ws.Range("J122").Value = 55558560
$ws.Range("L122").Value = 500027040
$ws.Range("N122").Value = -500031940
$ws.Range("H125").Value = 1843.3334
$ws.Range("I125").Value = 1180
$ws.Range("J125").Value = 2175
$ws.Range("K125").Value = 3540
$ws.Range("L125").Value = 6525
$ws.Range("M125").Value = 1380
$ws.Range("N125").Value = -16365
$ws.Range("H126").Value = 6976.6665
$ws.Range("I126").Value = 930
$ws.Range("K126").Value = 2790
$ws.Range("M126").Value = 2150
$ws.Range("H131").Value = 3621.9722
$ws.Range("I131").Value = 5349.1
$ws.Range("J131").Value = 2957.6924
$ws.Range("K131").Value = 16047.3
$ws.Range("L131").Value = 8873.0772
$ws.Range("M131").Value = -11007.3
$ws.Range("N131").Value = -18953.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2779
$ws.Range("I7").Value = 2735.8
$ws.Range("J7").Value = 2995
$ws.Range("K7").Value = 2735.8
$ws.Range("L7").Value = 2995
$ws.Range("M7").Value = -2623.8
$ws.Range("N7").Value = -3219
$ws.Range("H62").Value = 32500
$ws.Range("J62").Value = 32500
$ws.Range("L62").Value = 32500
$ws.Range("N62").Value = -33748
$ws.Range("H65").Value = 32500
$ws.Range("J65").Value = 32500
$ws.Range("L65").Value = 97500
$ws.Range("N65").Value = -103740
$ws.Range("H68").Value = 1727.3125
$ws.Range("I68").Value = 1634.96
$ws.Range("J68").Value = 2057.1428
$ws.Range("K68").Value = 1634.96
$ws.Range("L68").Value = 2057.1428
$ws.Range("M68").Value = -885.96
$ws.Range("N68").Value = -3555.1428
$ws.Range("H71").Value = 1727.3125
$ws.Range("I71").Value = 1634.96
$ws.Range("J71").Value = 2057.1428
$ws.Range("K71").Value = 8174.8
$ws.Range("L71").Value = 10285.714
$ws.Range("M71").Value = -4430.8
$ws.Range("N71").Value = -17773.714
$ws.Range("H126").Value = 2779
$ws.Range("I126").Value = 2735.8
$ws.Range("J126").Value = 2995
$ws.Range("K126").Value = 8207.400000000001
$ws.Range("L126").Value = 8985
$ws.Range("M126").Value = -5737.400000000001
$ws.Range("N126").Value = -13925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16367.842
$ws.Range("I62").Value = 16665.555
$ws.Range("J62").Value = 16099.9
$ws.Range("K62").Value = 16665.555
$ws.Range("L62").Value = 16099.9
$ws.Range("M62").Value = -16041.555
$ws.Range("N62").Value = -17347.9
$ws.Range("H65").Value = 16367.842
$ws.Range("I65").Value = 16665.555
$ws.Range("J65").Value = 16099.9
$ws.Range("K65").Value = 83327.77499999999
$ws.Range("L65").Value = 80499.5
$ws.Range("M65").Value = -80207.77499999999
$ws.Range("N65").Value = -86739.5
$ws.Range("H136").Value = 8149.8213
$ws.Range("I136").Value = 10313.667
$ws.Range("K136").Value = 30941.001
$ws.Range("M136").Value = -28391.001

